# Add ontology IDs to tags
# - Rename the "Computational Analysis" tag header (F13) to "Computation"
# - Add the accompanying Term Accession Number (F14) and Term Source REF (F15)
#   for the new "Computation" tag column on the isa_template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

$ws.Range("F13").Value = "Computation"
$ws.Range("F14").Value = "http://purl.obolibrary.org/obo/NCIT_C61298"
$ws.Range("F15").Value = "NCIT"
